# Update the "Förändrad" date column (C) for rows 2-14.
# The stored serial date value moves from 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180  # Column C ("Förändrad")
}
